$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest Adafruit IO reading (row 10) repeats the previous one (row 9) --
# timestamp 2024-09-25T18:06:40Z, temperature 25, N/A lat/long/elevation.
# Copy row 9 down so values/types line up exactly (keeps "25" as text
# instead of turning it into a number).
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial()
$excel.CutCopyMode = $false
